$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Keywords_Tag"
$ws.Range("A4").Value = "Update_Modification Date"
$ws.Range("A5").Value = "Theme_Category : Domain"

$ws.Range("D5").Select()
